# Fruta / hortaliza, semanal
#
# A new daily price record was inserted as row 17 (Fecha 19-08-2022,
# Volumen 15), pushing the existing rows 17-63 down to rows 18-64.
# The data for the new row 17 reuses the price/unit/origin figures that
# the (now) row 18 carries, with its own Fecha and Volumen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17; everything below
# (old rows 17-63) shifts down to rows 18-64, and the sheet's used range
# grows from A1:T63 to A1:T64 automatically.
$ws.Rows.Item(17).Insert()

# Seed the new row 17 by copying the row that used to be 17 (now sitting
# at row 18, right below), which carries the same Mercado/Producto
# metadata and the price/unit/origin figures the new record shares.
$ws.Range("A18:T18").Copy($ws.Range("A17:T17"))

# Overwrite just the two fields that differ for this new record.
$ws.Range("D17").Value = 44797   # Fecha: 19-08-2022
$ws.Range("M17").Value = 15      # Volumen
